# Update of archetype files:
# Add a new "ext_blind_test" material row to the "properties" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

$ws.Range("A12").Value = "ext_blind_test"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = "ignore"

# Match the printed page orientation recorded for both sheets.
$ws.PageSetup.Orientation = 1
$wb.Worksheets.Item("Sheet1").PageSetup.Orientation = 1

# Leave the cursor where the author last left it.
[void]$ws.Range("B15").Select()
